# edit.ps1 - applies the commit's changes to the presentation:
#   1. Updates the "datetimeFigureOut" date field text from 7/19/2018 to
#      7/20/2018 everywhere it appears (the slide master and every slide
#      layout in the deck).
#   2. On slide 1's subtitle, splits the run " and Experienced" into three
#      runs: " ", " and " and "Experienced" (matching the target edit,
#      which leaves the visible text as "Fresheres  and Experienced").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Part 1: date placeholder ("7/19/2018" -> "7/20/2018")
# ---------------------------------------------------------------------
$oldDate = "7/19/2018"
$newDate = "7/20/2018"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# Part 2: split the " and Experienced" run on slide 1's subtitle
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -like "*Fresheres*") {
            $target = $sh
        }
    }
}

if ($target -ne $null) {
    $tr = $target.TextFrame.TextRange
    $full = $tr.Text

    $needle = " and Experienced"
    $idx0 = $full.IndexOf($needle)

    if ($idx0 -ge 0) {
        $startPos = $idx0 + 1
        $needleLen = $needle.Length

        $expWord = "Experienced"
        $expLen = $expWord.Length
        $expStart = $startPos + ($needleLen - $expLen)

        $midLen = $needleLen - 1 - $expLen

        # Split "Experienced" off into its own run (same-length rewrite).
        $tr.Characters($expStart, $expLen).Text = $expWord

        # Split the leading single space into its own run.
        $tr.Characters($startPos, 1).Text = " "

        # Grow "and " into " and " (adds the extra leading space), which
        # becomes its own run distinct from the lone leading space run.
        $tr.Characters($startPos + 1, $midLen).Text = " and "
    }
}
